$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Top-of-table summary cells (rows 1-3 become "0M"; a couple of
# surrounding stat cells get refreshed numbers too).
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "1201"
$t.Cell(5,1).Range.Text  = "0.00001"
$t.Cell(6,1).Range.Text  = "0.00064"
$t.Cell(7,1).Range.Text  = "0.00016"
$t.Cell(9,1).Range.Text  = "0.00029"
$t.Cell(10,1).Range.Text = "0.00034"
$t.Cell(11,1).Range.Text = "0.00036"
$t.Cell(12,1).Range.Text = "0.21525"

# The last three rows previously held a whole tab-separated summary
# line each; they collapse down to the single lead figure (which now
# mirrors what used to sit at the very top of the table).
$t.Cell(44,1).Range.Text = "99.94"
$t.Cell(45,1).Range.Text = "0.22"
$t.Cell(46,1).Range.Text = "347"
